# The data rows (symbol, reel1..reel5) in the sheet were reshuffled into a
# different row order. Re-write the A:F values for the affected rows
# (3-18 and 21-23) to match the new order; everything else is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    3  = @(301, 6, 45, 30, 60, 45)
    4  = @(1202, 2, 10, 10, 10, 10)
    5  = @(1203, 3, 15, 15, 15, 15)
    6  = @(901, 16, 15, 45, 60, 60)
    7  = @(902, 1, 0, 0, 0, 0)
    8  = @(401, 9, 48, 67, 75, 45)
    9  = @(201, 9, 30, 15, 45, 30)
    10 = @(801, 3, 67, 65, 52, 45)
    11 = @(1201, 2, 10, 10, 10, 10)
    12 = @(501, 9, 52, 30, 75, 45)
    13 = @(601, 9, 60, 67, 60, 42)
    14 = @(1001, 18, 30, 75, 60, 72)
    15 = @(701, 3, 90, 45, 97, 15)
    16 = @(1101, 0, 15, 30, 30, 0)
    17 = @(2, 0, 2, 2, 2, 2)
    18 = @(3, 0, 3, 3, 3, 3)
    21 = @(1, 0, 2, 2, 2, 2)
    22 = @(402, 0, 0, 4, 0, 0)
    23 = @(602, 0, 0, 4, 0, 9)
}

$cols = @("A", "B", "C", "D", "E", "F")

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}
